# Rename the "congenital" category to "misc_long_term" across every
# variables_* worksheet in the workbook (new datasets + baseline regression).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Cells.Find("congenital")
    while ($cell) {
        $cell.Value = "misc_long_term"
        $cell = $ws.Cells.Find("congenital")
    }
}
